# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update the "K" column (column G) values that changed when the save data
# was regenerated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 0
    4  = 2
    5  = 0
    6  = 0
    7  = 2
    8  = 4
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 2
    14 = 0
    15 = 3
    16 = 1
    18 = 2
    19 = 1
    20 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
